# Adds a new "move zeros to left" data-structure exercise below the
# existing scratch area on the "Learning" sheet (rows 35-41), matching the
# commit "Added new forms and added new data structure problem that
# arrange all zero to left side of an given array".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 35: bold header row 0..5
$ws.Range("B35").Value = 0
$ws.Range("C35").Value = 1
$ws.Range("D35").Value = 2
$ws.Range("E35").Value = 3
$ws.Range("F35").Value = 4
$ws.Range("G35").Value = 5
$ws.Range("B35:G35").Font.Bold = $true

# Row 36: first pass of the algorithm
$ws.Range("B36").Value = 1
$ws.Range("C36").Value = 0
$ws.Range("D36").Value = 2
$ws.Range("E36").Value = 0
$ws.Range("F36").Value = 3
$ws.Range("G36").Value = 0

# Row 37: second pass
$ws.Range("B37").Value = 1
$ws.Range("C37").Value = 2
$ws.Range("D37").Value = 0
$ws.Range("E37").Value = 0
$ws.Range("F37").Value = 3
$ws.Range("G37").Value = 0

# Row 38: third pass / final arrangement
$ws.Range("B38").Value = 1
$ws.Range("C38").Value = 2
$ws.Range("D38").Value = 3
$ws.Range("E38").Value = 0
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 0

# Row 40: column labels for the index trace
$ws.Range("G40").Value = "I"
$ws.Range("H40").Value = "I+1"

# Row 41: trailing label
$ws.Range("E41").Value = "L"

# Leave the selection where the author left off editing.
$ws.Range("E41").Select()

# Best-effort: scroll the view so row 13 is at the top (matches the
# author's saved sheetView topLeftCell="A13").
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
